$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.633.98'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '1.642.26'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.77%  '
$ws.Range("E6").Value = '  +1.31%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +0.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0842'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("D12").Value = '1.870.94'
$ws.Range("E12").Value = '  +0.54%  '
$ws.Range("D13").Value = '1.660.53'
$ws.Range("E13").Value = '  +1.88%  '
$ws.Range("E14").Value = '  +2.46%  '
$ws.Range("E15").Value = '  +1.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.37'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.83%  '
$ws.Range("D17").Value = '26.672.78'
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '0.0₃0746'
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '216.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.86%  '
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("E21").Value = '  +1.32%  '
$ws.Range("E22").Value = '  +2.60%  '
$ws.Range("E23").Value = '  +1.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +12.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.69%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("E27").Value = '  -0.99%  '
$ws.Range("E28").Value = '  +4.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.77'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.49%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0518'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.59%  '
$ws.Range("E31").Value = '  +0.78%  '
$ws.Range("E32").Value = '  +2.27%  '
$ws.Range("E33").Value = '  +1.88%  '
$ws.Range("D34").Value = '1.277.41'
$ws.Range("E34").Value = '  +4.26%  '
$ws.Range("E35").Value = '  +2.82%  '
$ws.Range("E36").Value = '  +5.26%  '
$ws.Range("E37").Value = '  +0.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.534'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.69%  '
$ws.Range("E39").Value = '  +2.54%  '
$ws.Range("E40").Value = '  +0.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.817'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.70%  '
$ws.Range("E42").Value = '  -1.62%  '
$ws.Range("E43").Value = '  +2.07%  '
$ws.Range("D44").Value = '1.781.42'
$ws.Range("E44").Value = '  +0.74%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.06'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.85'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.19%  '
$ws.Range("E47").Value = '  +2.01%  '
$ws.Range("E48").Value = '  +0.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.83'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0969'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.78%  '
$ws.Range("E51").Value = '  -0.54%  '
